# Corrects the Sheet1 data table: fixes the mislabeled Q4-2021 row (A13 was a
# duplicate "IIIT_2021" label; it should be "IVT_2021"), and refreshes the
# Motociclista/Pasajero/Conductor counts (columns B, D, E) for rows 4-23 to
# their corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the row 13 label (was duplicated "IIIT_2021", should be "IVT_2021") ---
$ws.Cells.Item(13, 1).Value = "IVT_2021"

# --- Row 4 (IIIT_2019) ---
$ws.Cells.Item(4, 2).Value = 16
$ws.Cells.Item(4, 4).Value = 14
$ws.Cells.Item(4, 5).Value = 19

# --- Row 5 (IVT_2019) ---
$ws.Cells.Item(5, 4).Value = 15
$ws.Cells.Item(5, 5).Value = 11

# --- Row 6 (IT_2020) ---
$ws.Cells.Item(6, 5).Value = 9

# --- Row 7 (IIT_2020) ---
$ws.Cells.Item(7, 2).Value = 33
$ws.Cells.Item(7, 4).Value = 12
$ws.Cells.Item(7, 5).Value = 17

# --- Row 8 (IIIT_2020) ---
$ws.Cells.Item(8, 2).Value = 34
$ws.Cells.Item(8, 4).Value = 17

# --- Row 9 (IVT_2020) ---
$ws.Cells.Item(9, 2).Value = 42
$ws.Cells.Item(9, 4).Value = 17
$ws.Cells.Item(9, 5).Value = 16

# --- Row 10 (IT_2021) ---
$ws.Cells.Item(10, 2).Value = 46
$ws.Cells.Item(10, 4).Value = 20
$ws.Cells.Item(10, 5).Value = 18

# --- Row 11 (IIT_2021) ---
$ws.Cells.Item(11, 2).Value = 28
$ws.Cells.Item(11, 4).Value = 12
$ws.Cells.Item(11, 5).Value = 16

# --- Row 12 (IIIT_2021) ---
$ws.Cells.Item(12, 2).Value = 33
$ws.Cells.Item(12, 4).Value = 15

# --- Row 13 (IVT_2021) ---
$ws.Cells.Item(13, 2).Value = 34
$ws.Cells.Item(13, 4).Value = 13
$ws.Cells.Item(13, 5).Value = 14
$ws.Cells.Item(13, 6).Value = 4

# --- Row 14 (IT_2022) ---
$ws.Cells.Item(14, 2).Value = 43
$ws.Cells.Item(14, 4).Value = 17
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 9

# --- Row 15 (IIT_2022) ---
$ws.Cells.Item(15, 2).Value = 47
$ws.Cells.Item(15, 4).Value = 24
$ws.Cells.Item(15, 5).Value = 11
$ws.Cells.Item(15, 6).Value = 3

# --- Row 16 (IIIT_2022) ---
$ws.Cells.Item(16, 3).Value = 48
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 12

# --- Row 17 (IVT_2022) ---
$ws.Cells.Item(17, 3).Value = 46
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = 4

# --- Row 18 (IT_2023) ---
$ws.Cells.Item(18, 2).Value = 33
$ws.Cells.Item(18, 4).Value = 14
$ws.Cells.Item(18, 5).Value = 7

# --- Row 19 (IIT_2023) ---
$ws.Cells.Item(19, 2).Value = 43
$ws.Cells.Item(19, 4).Value = 20
$ws.Cells.Item(19, 5).Value = 9

# --- Row 20 (IIIT_2023) ---
$ws.Cells.Item(20, 2).Value = 58
$ws.Cells.Item(20, 4).Value = 12
$ws.Cells.Item(20, 5).Value = 21

# --- Row 21 (IVT_2023) ---
$ws.Cells.Item(21, 2).Value = 70
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 10

# --- Row 22 (IT_2024) ---
$ws.Cells.Item(22, 2).Value = 55
$ws.Cells.Item(22, 4).Value = 9
$ws.Cells.Item(22, 5).Value = 12

# --- Row 23 (IIT_2024) ---
$ws.Cells.Item(23, 2).Value = 63
